# Reorders the per-record data (Fecha, Volumen, Precio mínimo/máximo/promedio,
# Origen, Precio $/Kg) across data rows 2-19 of the sheet, as described by the
# commit "Fruta / hortaliza, semanal". The descriptive columns (Mercado ID,
# Mercado, Región, Codreg, Tipo, Producto*, Categoría*, Variedad, Calidad,
# Unidad de comercialización, Kg/unidad) are identical on every row, so the
# net effect of the edit is a permutation of whole data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Maps destination row number -> source row number (i.e. row $key ends up
# holding the data that used to live in row $map[$key]).
$map = @{
    2  = 10
    3  = 14
    4  = 7
    5  = 9
    6  = 19
    7  = 18
    8  = 11
    9  = 6
    10 = 2
    11 = 3
    12 = 4
    13 = 13
    14 = 5
    15 = 17
    16 = 12
    17 = 8
    18 = 16
    19 = 15
}

$lastCol = "T"

# Snapshot every source row's full contents (A:T) before writing anything,
# so that later writes never clobber data that's still needed as a source.
$snapshot = @{}
foreach ($r in $map.Keys) {
    $snapshot[$r] = $ws.Range("A$r`:$lastCol$r").Value2
}

foreach ($destRow in $map.Keys) {
    $srcRow = $map[$destRow]
    $ws.Range("A$destRow`:$lastCol$destRow").Value = $snapshot[$srcRow]
}
